# edit.ps1
# Applies the changes described by the diff:
#  1. Merge split runs (re-saved by Word as a single run with identical rPr)
#     in three places:
#       - "...to identify predefined vulnerabilities..." paragraph
#       - "IT – Information Technology" glossary entry
#       - "CLI – Command Line Interface" glossary entry
#  2. Center the page number in the default (primary) footer.
#  3. Turn off automatic hyphenation for the Normal style
#     (<w:suppressAutoHyphens/> in styles.xml).

$d = $word.ActiveDocument

# --- 1a. Merge the split runs in the Metasploit / predefined-vulnerabilities paragraph ---
$text1 = "Moving beyond assessment, our plan is to use exploitation frameworks such as Metasploit to actively test the system. We will use the previous reconnaissance methodologies to identify predefined vulnerabilities and match them based on a curated list. This last step will emphasize how vulnerable their system is."
$found1 = $d.Content.Find.Execute($text1, $true, $false, $false, $false, $false, $true, 1, $false, $text1, 2)

# --- 1b. Merge "IT – " + "Information Technology" into one run ---
$text2 = "IT – Information Technology"
$found2 = $d.Content.Find.Execute($text2, $true, $false, $false, $false, $false, $true, 1, $false, $text2, 2)

# --- 1c. Merge "C" + "LI – Command Line Interface" into one run ---
$text3 = "CLI – Command Line Interface"
$found3 = $d.Content.Find.Execute($text3, $true, $false, $false, $false, $false, $true, 1, $false, $text3, 2)

# --- 2. Center the page number in the default (primary) footer ---
$footer = $d.Sections(1).Footers(1)
$footer.Range.ParagraphFormat.Alignment = 1

# --- 3. Suppress automatic hyphenation on the Normal style ---
$normalStyle = $d.Styles("Normal")
$normalStyle.ParagraphFormat.Hyphenation = $false

Write-Output "merge1=$found1 merge2=$found2 merge3=$found3"
